# Fruta / hortaliza, semanal
#
# A new weekly record is inserted at row 203 (pushing the existing rows
# 203..231 down to 204..232, same as the source OOXML diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 203, shifting everything below it down by one.
$ws.Rows.Item(203).Insert()

# Populate the new row with the latest weekly observation.
$ws.Cells.Item(203, 1).Value = 3
$ws.Cells.Item(203, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(203, 3).Value = "Coquimbo"
$ws.Cells.Item(203, 4).Value = 44491
$ws.Cells.Item(203, 5).Value = 5
$ws.Cells.Item(203, 6).Value = 100112031
$ws.Cells.Item(203, 7).Value = "Poroto verde"
$ws.Cells.Item(203, 8).Value = "Magnum"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 38
$ws.Cells.Item(203, 11).Value = 42000
$ws.Cells.Item(203, 12).Value = 42000
$ws.Cells.Item(203, 13).Value = 42000
$ws.Cells.Item(203, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(203, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(203, 16).Value = 1680
$ws.Cells.Item(203, 17).Value = 25
$ws.Cells.Item(203, 18).Value = "Hortaliza"

# Keep the date cell formatted the same way as the rest of column D.
$ws.Cells.Item(203, 4).NumberFormat = $ws.Cells.Item(204, 4).NumberFormat
